$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells for the team record columns
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Copy the style of an existing header cell (e.g. AC1) onto the new headers
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = 0

# Fill in the team record values (Wins=61, Losses=101, Ties=0) for every data row
$lastRow = 66
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 30).Value = 61   # AD
    $ws.Cells.Item($r, 31).Value = 101  # AE
    $ws.Cells.Item($r, 32).Value = 0    # AF
}
